{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the LOM3219.docx content revision:\n//   - bump activation date to 2025\n//   - replace PT/EN \"Objetivos\" paragraphs\n//   - replace PT/EN \"Programa\" paragraphs\n//   - append merged Crit\u00e9rio/Norma de Recupera\u00e7\u00e3o text onto the\n//     \"Aulas expositivas e semin\u00e1rios.\" run\n//   - replace the Bibliografia paragraph text\n//\n// Each change is applied as an exact, whole-string search-and-replace so\n// that only the targeted run's text is touched; surrounding runs/formatting\n// (bold labels, line breaks, italics, etc.) stay untouched.\n\nconst replacements = [\n  [\n    \"Ativa\u00e7\u00e3o: 01/01/2023\",\n    \"Ativa\u00e7\u00e3o: 01/01/2025\",\n  ],\n  [\n    \"Apresentar os conceitos de nanoci\u00eancia e nanotecnologia. As propriedades f\u00edsicas e qu\u00edmicas dos materiais em escala nanom\u00e9trica s\u00e3o descritas pelas leis da mec\u00e2nica qu\u00e2ntica, apresentando nessas dimens\u00f5es caracter\u00edsticas distintas dos materiais em escala macrosc\u00f3pica. O conhecimento dessa \u00e1rea interdisciplinar \u00e9 fundamental na forma\u00e7\u00e3o de um pesquisador e/ou um profissional atuando na \u00e1rea de materiais.\",\n    \"A matura\u00e7\u00e3o da nanotecnologia revelou que se trata de uma disciplina \u00fanica e distinta, em vez de uma especializa\u00e7\u00e3o dentro de um campo maior. Um curso sobre esse assunto envolve qu\u00edmica, f\u00edsica e engenharia focada em Nano. Deve ser integrado, multidisciplinar e especificamente em Nano. A ideia \u00e9 construir uma base s\u00f3lida nos m\u00e9todos de caracteriza\u00e7\u00e3o e fabrica\u00e7\u00e3o enquanto integra a f\u00edsicas e a qu\u00edmica relevantes aos problemas envolvidos. Examinando os aspectos de engenharia, bem como nanomateriais e aplica\u00e7\u00f5es espec\u00edficas nos setores de energia e eletr\u00f4nica.\",\n  ],\n  [\n    \"Present the concepts of nanoscience and nanotechnology. The physical and chemical properties of materials on a nanometer scale are described by the laws of quantum mechanics, presenting in these dimensions different characteristics of materials on a macroscopic scale. The knowledge of this interdisciplinary area is fundamental in the formation of a researcher and/or a professional working in the area of materials.\",\n    \"Nanotechnology maturation has revealed that it is a unique and distinct discipline rather than a specialization within a larger field. A course on this subject involves chemistry, physics and engineering focused on Nano. It must be integrated, multidisciplinary and specifically in nano. The idea is to build a solid foundation on characterization and manufacturing methods while integrating with physical and chemistry relevant to the problems involved. Examining engineering aspects as well as nanomaterials and specific applications in the energy and electronics sectors.\",\n  ],\n  [\n    \"Conceitua\u00e7\u00e3o: nanoci\u00eancia e nanotecnologia. Sistemas de baixa dimensionalidade. Confinamento qu\u00e2ntico. Liga\u00e7\u00f5es qu\u00edmicas: mol\u00e9culas e aglomerados. Propriedades eletr\u00f4nicas e estruturais. S\u00edntese e fabrica\u00e7\u00e3o de materiais em escala nanom\u00e9trica: t\u00e9cnicas de baixo para cima (bottom-up) e de cima para baixo (top-down). Fullerenos e nanotubos de carbono. Autoorganiza\u00e7\u00e3o molecular e sistemas supramoleculares. Fios e pontos qu\u00e2nticos. Nanopart\u00edculas magn\u00e9ticas. T\u00e9cnicas de caracteriza\u00e7\u00e3o: difra\u00e7\u00e3o, espalhamento e absor\u00e7\u00e3o de raios X, microscopia de varredura por tunelamento (STM), microscopia de for\u00e7a at\u00f4mica (AFM), microscopia eletr\u00f4nica de transmiss\u00e3o. Propriedades de transporte: transporte bal\u00edstico, condut\u00e2ncia qu\u00e2ntica, bloqueio coulombiano. Dispositivos moleculares. Transporte difusivo. Nanomagnetismo: ordem magn\u00e9tica, superparamagnetismo e Spintr\u00f4nica. Aplica\u00e7\u00f5es.\",\n    \"Perspectivas: nanoci\u00eancia e nanotecnologia - a distin\u00e7\u00e3o; Implica\u00e7\u00f5es sociais de nanoNanotools: m\u00e9todos de caracteriza\u00e7\u00e3o; M\u00e9todos de fabrica\u00e7\u00e3oF\u00edsica: Propriedades e fen\u00f4menos: materiais, estrutura e nanosurface; Energia na nanoescalaQu\u00edmica: s\u00edntese e modifica\u00e7\u00e3o: nanomateriais \u00e0 base de carbono; Intera\u00e7\u00f5es qu\u00edmicas na nanoescalaAplica\u00e7\u00f5es: nanoetronics; nanomagnetismo; nanomec\u00e2nica\",\n  ],\n  [\n    \"Conceptualization: nanoscience and nanotechnology. Low-dimensional systems. Quantum Confinement. Chemical bonds: molecules and clusters. Electronic and structural properties. Synthesis and fabrication of materials at the nanometer scale: bottom-up and top-down techniques. Fullerenes and carbon nanotubes. Molecular self-organization and supramolecular systems. Quantum wires and dots. Magnetic nanoparticles. Characterization techniques: X-ray diffraction, scattering and absorption, scanning tunneling microscopy (STM), atomic force microscopy (AFM), transmission electron microscopy. Transport properties: ballistic transport, quantum conductance, Coulomb blocking. Molecular devices. Diffusive transport. Nanomagnetism: magnetic order, superparamagnetism and spintronics. Applications.\",\n    \"Perspectives: Nanoscience and Nanotechnology\u2014The Distinction; Societal Implications of NanoNanotools: Characterization Methods; Fabrication MethodsPhysics: Properties and Phenomena: Materials, Structure, and the Nanosurface; Energy at the NanoscaleChemistry: Synthesis and Modification: Carbon-Based Nanomaterials; Chemical Interactions at the NanoscaleApplications: nanoeletronics; nanomagnetism; nanomechanics\",\n  ],\n  [\n    \"Aulas expositivas e semin\u00e1rios.\",\n    \"Aulas expositivas e semin\u00e1rios.Crit\u00e9rioDuas provas escritas: conceitos P1 e P2. Conceito Final = (P1 + 2P2)/3Norma de Recupera\u00e7\u00e3oAplica\u00e7\u00e3o de uma prova escrita dentro do prazo regimental antes do in\u00edcio do pr\u00f3ximo semestre letivo. A nota da segunda avalia\u00e7\u00e3o ser\u00e1 a m\u00e9dia aritm\u00e9tica entre a nota da prova de recupera\u00e7\u00e3o e a nota final da primeira avalia\u00e7\u00e3o\",\n  ],\n  [\n    \"TIMP, G. Nanotechnology, Springer, 1998. FERRY, D. K. Transport in Nanostructures, Cambridge University Press, 1999. WASER, R. Nanoelectronics and Information Technology, Wiley-UCM, 2003. DATTA, S. Quantum Transport: Atom to Transistor, Cambridge University Press, 2005. RATNER, M.; RATNER, D. Nanotechnology, Prentice Hall, 2003. DRESSELHAUS, M. Physical Properties of Carbon Nanotubes, Imperial College Press, 1998.\",\n    \"Gabor L. Hornyak, H.F. Tibbals, Joydeep Dutta, John J. Moore. Introduction to Nanoscience and Nanotechnology. CRC Press. 2009TIMP, G. Nanotechnology, Springer, 1998.Bhushan, B. (ed.) Springer Handbook of Nanotechnology, Springer, 2010.\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"Could not find expected text: \" + oldText.substring(0, 60));\n  }\n\n  // Replace the matched range's text in place; this preserves the run's\n  // existing formatting (bold/italic/etc.) and keeps sibling runs (e.g. the\n  // bold \"Crit\u00e9rio: \" label run) untouched.\n  found.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the LOM3219.docx content revision:\n#   - bump activation date to 2025\n#   - replace PT/EN \"Objetivos\" paragraphs\n#   - replace PT/EN \"Programa\" paragraphs\n#   - append merged Criterio/Norma de Recuperacao text onto the\n#     \"Aulas expositivas e seminarios.\" run\n#   - replace the Bibliografia paragraph text\n#\n# Each change is a literal (non-wildcard) Find/Replace over the whole\n# document body, scoped to an exact, unique source string so only the\n# targeted run's text changes; sibling runs (e.g. the bold \"Criterio: \"\n# label) and their formatting are left alone.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $ok) {\n        throw \"Could not find expected text: $findText\"\n    }\n}\n\nReplace-ExactText \"Ativa\u00e7\u00e3o: 01/01/2023\" \"Ativa\u00e7\u00e3o: 01/01/2025\"\n\nReplace-ExactText \"Apresentar os conceitos de nanoci\u00eancia e nanotecnologia. As propriedades f\u00edsicas e qu\u00edmicas dos materiais em escala nanom\u00e9trica s\u00e3o descritas pelas leis da mec\u00e2nica qu\u00e2ntica, apresentando nessas dimens\u00f5es caracter\u00edsticas distintas dos materiais em escala macrosc\u00f3pica. O conhecimento dessa \u00e1rea interdisciplinar \u00e9 fundamental na forma\u00e7\u00e3o de um pesquisador e/ou um profissional atuando na \u00e1rea de materiais.\" \"A matura\u00e7\u00e3o da nanotecnologia revelou que se trata de uma disciplina \u00fanica e distinta, em vez de uma especializa\u00e7\u00e3o dentro de um campo maior. Um curso sobre esse assunto envolve qu\u00edmica, f\u00edsica e engenharia focada em Nano. Deve ser integrado, multidisciplinar e especificamente em Nano. A ideia \u00e9 construir uma base s\u00f3lida nos m\u00e9todos de caracteriza\u00e7\u00e3o e fabrica\u00e7\u00e3o enquanto integra a f\u00edsicas e a qu\u00edmica relevantes aos problemas envolvidos. Examinando os aspectos de engenharia, bem como nanomateriais e aplica\u00e7\u00f5es espec\u00edficas nos setores de energia e eletr\u00f4nica.\"\n\nReplace-ExactText \"Present the concepts of nanoscience and nanotechnology. The physical and chemical properties of materials on a nanometer scale are described by the laws of quantum mechanics, presenting in these dimensions different characteristics of materials on a macroscopic scale. The knowledge of this interdisciplinary area is fundamental in the formation of a researcher and/or a professional working in the area of materials.\" \"Nanotechnology maturation has revealed that it is a unique and distinct discipline rather than a specialization within a larger field. A course on this subject involves chemistry, physics and engineering focused on Nano. It must be integrated, multidisciplinary and specifically in nano. The idea is to build a solid foundation on characterization and manufacturing methods while integrating with physical and chemistry relevant to the problems involved. Examining engineering aspects as well as nanomaterials and specific applications in the energy and electronics sectors.\"\n\nReplace-ExactText \"Conceitua\u00e7\u00e3o: nanoci\u00eancia e nanotecnologia. Sistemas de baixa dimensionalidade. Confinamento qu\u00e2ntico. Liga\u00e7\u00f5es qu\u00edmicas: mol\u00e9culas e aglomerados. Propriedades eletr\u00f4nicas e estruturais. S\u00edntese e fabrica\u00e7\u00e3o de materiais em escala nanom\u00e9trica: t\u00e9cnicas de baixo para cima (bottom-up) e de cima para baixo (top-down). Fullerenos e nanotubos de carbono. Autoorganiza\u00e7\u00e3o molecular e sistemas supramoleculares. Fios e pontos qu\u00e2nticos. Nanopart\u00edculas magn\u00e9ticas. T\u00e9cnicas de caracteriza\u00e7\u00e3o: difra\u00e7\u00e3o, espalhamento e absor\u00e7\u00e3o de raios X, microscopia de varredura por tunelamento (STM), microscopia de for\u00e7a at\u00f4mica (AFM), microscopia eletr\u00f4nica de transmiss\u00e3o. Propriedades de transporte: transporte bal\u00edstico, condut\u00e2ncia qu\u00e2ntica, bloqueio coulombiano. Dispositivos moleculares. Transporte difusivo. Nanomagnetismo: ordem magn\u00e9tica, superparamagnetismo e Spintr\u00f4nica. Aplica\u00e7\u00f5es.\" \"Perspectivas: nanoci\u00eancia e nanotecnologia - a distin\u00e7\u00e3o; Implica\u00e7\u00f5es sociais de nanoNanotools: m\u00e9todos de caracteriza\u00e7\u00e3o; M\u00e9todos de fabrica\u00e7\u00e3oF\u00edsica: Propriedades e fen\u00f4menos: materiais, estrutura e nanosurface; Energia na nanoescalaQu\u00edmica: s\u00edntese e modifica\u00e7\u00e3o: nanomateriais \u00e0 base de carbono; Intera\u00e7\u00f5es qu\u00edmicas na nanoescalaAplica\u00e7\u00f5es: nanoetronics; nanomagnetismo; nanomec\u00e2nica\"\n\nReplace-ExactText \"Conceptualization: nanoscience and nanotechnology. Low-dimensional systems. Quantum Confinement. Chemical bonds: molecules and clusters. Electronic and structural properties. Synthesis and fabrication of materials at the nanometer scale: bottom-up and top-down techniques. Fullerenes and carbon nanotubes. Molecular self-organization and supramolecular systems. Quantum wires and dots. Magnetic nanoparticles. Characterization techniques: X-ray diffraction, scattering and absorption, scanning tunneling microscopy (STM), atomic force microscopy (AFM), transmission electron microscopy. Transport properties: ballistic transport, quantum conductance, Coulomb blocking. Molecular devices. Diffusive transport. Nanomagnetism: magnetic order, superparamagnetism and spintronics. Applications.\" \"Perspectives: Nanoscience and Nanotechnology\u2014The Distinction; Societal Implications of NanoNanotools: Characterization Methods; Fabrication MethodsPhysics: Properties and Phenomena: Materials, Structure, and the Nanosurface; Energy at the NanoscaleChemistry: Synthesis and Modification: Carbon-Based Nanomaterials; Chemical Interactions at the NanoscaleApplications: nanoeletronics; nanomagnetism; nanomechanics\"\n\nReplace-ExactText \"Aulas expositivas e semin\u00e1rios.\" \"Aulas expositivas e semin\u00e1rios.Crit\u00e9rioDuas provas escritas: conceitos P1 e P2. Conceito Final = (P1 + 2P2)/3Norma de Recupera\u00e7\u00e3oAplica\u00e7\u00e3o de uma prova escrita dentro do prazo regimental antes do in\u00edcio do pr\u00f3ximo semestre letivo. A nota da segunda avalia\u00e7\u00e3o ser\u00e1 a m\u00e9dia aritm\u00e9tica entre a nota da prova de recupera\u00e7\u00e3o e a nota final da primeira avalia\u00e7\u00e3o\"\n\nReplace-ExactText \"TIMP, G. Nanotechnology, Springer, 1998. FERRY, D. K. Transport in Nanostructures, Cambridge University Press, 1999. WASER, R. Nanoelectronics and Information Technology, Wiley-UCM, 2003. DATTA, S. Quantum Transport: Atom to Transistor, Cambridge University Press, 2005. RATNER, M.; RATNER, D. Nanotechnology, Prentice Hall, 2003. DRESSELHAUS, M. Physical Properties of Carbon Nanotubes, Imperial College Press, 1998.\" \"Gabor L. Hornyak, H.F. Tibbals, Joydeep Dutta, John J. Moore. Introduction to Nanoscience and Nanotechnology. CRC Press. 2009TIMP, G. Nanotechnology, Springer, 1998.Bhushan, B. (ed.) Springer Handbook of Nanotechnology, Springer, 2010.\"\n"}
